$wb = $excel.ActiveWorkbook

# OFF sheet - Week 17 Home row update
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 241
$wsOff.Range("C2").Value = 162
$wsOff.Range("D2").Value = 55
$wsOff.Range("F2").Value = 5
$wsOff.Range("G2").Value = 2

# DEF sheet - Week 17 Home row update
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 213
$wsDef.Range("C2").Value = 131
$wsDef.Range("D2").Value = 38
$wsDef.Range("E2").Value = 11

$wb.Save()
